# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Update "last updated" timestamp in title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 13:05"

# 2) Estados Unidos (row 4): update totals / new cases / recovered
$ws.Range("B4").Value = 1593297
$ws.Range("C4").Value = 574
$ws.Range("E4").Value = 1127485

# 3) Iran (row 13): update totals / new cases / active / recovered / critical / deaths
$ws.Range("B13").Value = 129341
$ws.Range("C13").Value = 2392
$ws.Range("D13").Value = 100564
$ws.Range("E13").Value = 21528
$ws.Range("G13").Value = 66
$ws.Range("H13").Value = 7249

# 4) Suiza (row 28): update totals / new cases / recovered
$ws.Range("B28").Value = 30694
$ws.Range("C28").Value = 36
$ws.Range("E28").Value = 1002

# 5) Nepal overtakes Reunion in the ranking (sorted descending by total cases).
#    Row 131 now holds Nepal with its updated figures; row 132 now holds
#    Reunion with its previous (unchanged) figures.
$ws.Range("A131").Value = "Nepal"
$ws.Range("B131").Value = 453
$ws.Range("C131").Value = 26
$ws.Range("D131").Value = 49
$ws.Range("E131").Value = 401
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 3

$ws.Range("A132").Value = "Reunion"
$ws.Range("B132").Value = 447
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 411
$ws.Range("E132").Value = 35
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 1

$wb.Save()
